# Commit: "a new row added"
#
# The "Programming Basics Book - Plan" sheet gains one more row (row 14)
# appended right after the last existing data row (row 13): a new chapter
# item "#11" with name "Нов ред" in column B. Columns C..H are left blank,
# matching the source edit. The new row simply inherits the column default
# formatting (no explicit per-cell style is introduced), exactly like a
# value freshly typed into previously-empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row (A14/B14) right below the previous last row (13).
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Нов ред"

# Reflect the author having scrolled down and landed on the cell right
# below the freshly-typed row.
$ws.Range("B15").Select() | Out-Null
